$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text, matching the source data (inline strings).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.606.74'
$ws.Range("E2").Value = '  -3.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.893.26'
$ws.Range("E3").Value = '  -3.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.00'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.78'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.877.92'
$ws.Range("E7").Value = '  -3.94%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.675'
$ws.Range("E8").Value = '  -6.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.724'
$ws.Range("E10").Value = '  -4.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("E11").Value = '  -6.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.05'
$ws.Range("E12").Value = '  +14.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000310'
$ws.Range("E13").Value = '  -4.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.37'
$ws.Range("E14").Value = '  -4.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.485.72'
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.865.76'
$ws.Range("E16").Value = '  -4.91%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.12'
$ws.Range("E17").Value = '  -4.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.52'
$ws.Range("E18").Value = '  -4.03%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.131'
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.15'
$ws.Range("E20").Value = '  -4.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.354.20'
$ws.Range("E21").Value = '  -4.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '417.13'
$ws.Range("E22").Value = '  -5.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.73'
$ws.Range("E23").Value = '  -9.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.43'
$ws.Range("E24").Value = '  -4.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.08'
$ws.Range("E25").Value = '  +1.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.02'
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.08'
$ws.Range("E27").Value = '  -3.07%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.37'
$ws.Range("E28").Value = '  -5.92%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.82'
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.59'
$ws.Range("E30").Value = '  +15.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.73'
$ws.Range("E31").Value = '  -5.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.37'
$ws.Range("E32").Value = '  +8.17%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '48.28'
$ws.Range("E33").Value = '  +14.07%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '670.30'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.99'
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.125'
$ws.Range("E36").Value = '  -2.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '63.24'
$ws.Range("E37").Value = '  -6.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.423'
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.146'
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.36'
$ws.Range("E40").Value = '  -5.52%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0793'
$ws.Range("E41").Value = '  -8.01%  '
$ws.Range("B42").Value = 'Dai'
$ws.Range("C42").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.20'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0471'
$ws.Range("E45").Value = '  -4.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.145'
$ws.Range("E46").Value = '  -7.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.66'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.42'
$ws.Range("E48").Value = '  +4.21%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.30'
$ws.Range("E49").Value = '  -4.64%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000271'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.91'
$ws.Range("E51").Value = '  -5.54%  '
